$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn = $wb.Worksheets.Item(2)
$wsDeDe = $wb.Worksheets.Item(3)

# Update handoff/handback timestamp cells
$wsZhCn.Range('E2').Value2 = '2016-03-12 18:51:09'
$wsZhCn.Range('H2').Value2 = '2016-03-12 18:51:25'
$wsZhCn.Range('E3').Value2 = '2016-03-12 18:51:09'
$wsZhCn.Range('H3').Value2 = '2016-03-12 18:51:25'
$wsDeDe.Range('E2').Value2 = '2016-03-12 18:51:12'
$wsDeDe.Range('H2').Value2 = '2016-03-12 18:51:31'
$wsDeDe.Range('E3').Value2 = '2016-03-12 18:51:12'
$wsDeDe.Range('H3').Value2 = '2016-03-12 18:51:31'

# Update hyperlinked cells (source/target file names) - recreate hyperlink
# in place so the underlying address/relationship is preserved exactly,
# while only the visible display text changes.
$wsOverview.Range('A2').Hyperlinks.Delete()
$wsOverview.Range('A3').Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/6999aef4a2cfebfef8412ba27bb42cd5302f70b9/e2e/92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.md', "", "", '3a7cbadf-f743-4554-8d66-08dcba5cdbdb.md')
$wsOverview.Hyperlinks.Add($wsOverview.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/6999aef4a2cfebfef8412ba27bb42cd5302f70b9/e2e/d0dc96df-64c6-47fb-94d3-fbc50197c361.md', "", "", 'ffff485b18b2-a352-4f19-9ae5-30fcd29e5a8e.md')
$wsZhCn.Range('A2').Hyperlinks.Delete()
$wsZhCn.Range('D2').Hyperlinks.Delete()
$wsZhCn.Range('F2').Hyperlinks.Delete()
$wsZhCn.Range('G2').Hyperlinks.Delete()
$wsZhCn.Range('A3').Hyperlinks.Delete()
$wsZhCn.Range('D3').Hyperlinks.Delete()
$wsZhCn.Range('F3').Hyperlinks.Delete()
$wsZhCn.Range('G3').Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/6999aef4a2cfebfef8412ba27bb42cd5302f70b9/e2e/92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.md', "", "", '3a7cbadf-f743-4554-8d66-08dcba5cdbdb.md')
$wsZhCn.Hyperlinks.Add($wsZhCn.Range('D2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2258437a52de4fa8b7da741ad685c297c204e0f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.1a7532fc668ee458d29b28e4d0235919d447ef1d.zh-cn.xlf', "", "", '3a7cbadf-f743-4554-8d66-08dcba5cdbdb.f7024155cf44173efde75242c5a93109592061fd.zh-cn.xlf')
$wsZhCn.Hyperlinks.Add($wsZhCn.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/33f14ea79b0100d329cc4f0629971e5ae3e4135e/e2e/92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.md', "", "", '3a7cbadf-f743-4554-8d66-08dcba5cdbdb.md')
$wsZhCn.Hyperlinks.Add($wsZhCn.Range('G2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/d4644741f8b86468e374b0a4d807a5f5a0c0878f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.1a7532fc668ee458d29b28e4d0235919d447ef1d.zh-cn.xlf', "", "", '3a7cbadf-f743-4554-8d66-08dcba5cdbdb.f7024155cf44173efde75242c5a93109592061fd.zh-cn.xlf')
$wsZhCn.Hyperlinks.Add($wsZhCn.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/6999aef4a2cfebfef8412ba27bb42cd5302f70b9/e2e/d0dc96df-64c6-47fb-94d3-fbc50197c361.md', "", "", 'ffff485b18b2-a352-4f19-9ae5-30fcd29e5a8e.md')
$wsZhCn.Hyperlinks.Add($wsZhCn.Range('D3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2258437a52de4fa8b7da741ad685c297c204e0f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d0dc96df-64c6-47fb-94d3-fbc50197c361.22ad50810204c9096ca26f259094112d8ab2d6ff.zh-cn.xlf', "", "", '3a7cbadf-f743-4554-8d66-08dcba5cdbdb.f7024155cf44173efde75242c5a93109592061fd.zh-cn.xlf')
$wsZhCn.Hyperlinks.Add($wsZhCn.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/33f14ea79b0100d329cc4f0629971e5ae3e4135e/e2e/d0dc96df-64c6-47fb-94d3-fbc50197c361.md', "", "", 'ffff485b18b2-a352-4f19-9ae5-30fcd29e5a8e.md')
$wsZhCn.Hyperlinks.Add($wsZhCn.Range('G3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/d4644741f8b86468e374b0a4d807a5f5a0c0878f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d0dc96df-64c6-47fb-94d3-fbc50197c361.22ad50810204c9096ca26f259094112d8ab2d6ff.zh-cn.xlf', "", "", '3a7cbadf-f743-4554-8d66-08dcba5cdbdb.f7024155cf44173efde75242c5a93109592061fd.zh-cn.xlf')
$wsDeDe.Range('A2').Hyperlinks.Delete()
$wsDeDe.Range('D2').Hyperlinks.Delete()
$wsDeDe.Range('F2').Hyperlinks.Delete()
$wsDeDe.Range('G2').Hyperlinks.Delete()
$wsDeDe.Range('A3').Hyperlinks.Delete()
$wsDeDe.Range('D3').Hyperlinks.Delete()
$wsDeDe.Range('F3').Hyperlinks.Delete()
$wsDeDe.Range('G3').Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/6999aef4a2cfebfef8412ba27bb42cd5302f70b9/e2e/92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.md', "", "", '3a7cbadf-f743-4554-8d66-08dcba5cdbdb.md')
$wsDeDe.Hyperlinks.Add($wsDeDe.Range('D2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc84705b8a1ef2a402362b46e0fbd6e9edeec32d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.1a7532fc668ee458d29b28e4d0235919d447ef1d.de-de.xlf', "", "", '3a7cbadf-f743-4554-8d66-08dcba5cdbdb.f7024155cf44173efde75242c5a93109592061fd.de-de.xlf')
$wsDeDe.Hyperlinks.Add($wsDeDe.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/76a4bf33a238ba398cb5740c7dda43475c54a83c/e2e/92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.md', "", "", '3a7cbadf-f743-4554-8d66-08dcba5cdbdb.md')
$wsDeDe.Hyperlinks.Add($wsDeDe.Range('G2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/449dc07c859b62d0376b06b48dbe0017d2c350c5/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.1a7532fc668ee458d29b28e4d0235919d447ef1d.de-de.xlf', "", "", '3a7cbadf-f743-4554-8d66-08dcba5cdbdb.f7024155cf44173efde75242c5a93109592061fd.de-de.xlf')
$wsDeDe.Hyperlinks.Add($wsDeDe.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/6999aef4a2cfebfef8412ba27bb42cd5302f70b9/e2e/d0dc96df-64c6-47fb-94d3-fbc50197c361.md', "", "", 'ffff485b18b2-a352-4f19-9ae5-30fcd29e5a8e.md')
$wsDeDe.Hyperlinks.Add($wsDeDe.Range('D3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc84705b8a1ef2a402362b46e0fbd6e9edeec32d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d0dc96df-64c6-47fb-94d3-fbc50197c361.22ad50810204c9096ca26f259094112d8ab2d6ff.de-de.xlf', "", "", '3a7cbadf-f743-4554-8d66-08dcba5cdbdb.f7024155cf44173efde75242c5a93109592061fd.de-de.xlf')
$wsDeDe.Hyperlinks.Add($wsDeDe.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/76a4bf33a238ba398cb5740c7dda43475c54a83c/e2e/d0dc96df-64c6-47fb-94d3-fbc50197c361.md', "", "", 'ffff485b18b2-a352-4f19-9ae5-30fcd29e5a8e.md')
$wsDeDe.Hyperlinks.Add($wsDeDe.Range('G3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/449dc07c859b62d0376b06b48dbe0017d2c350c5/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d0dc96df-64c6-47fb-94d3-fbc50197c361.22ad50810204c9096ca26f259094112d8ab2d6ff.de-de.xlf', "", "", '3a7cbadf-f743-4554-8d66-08dcba5cdbdb.f7024155cf44173efde75242c5a93109592061fd.de-de.xlf')

